$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right and Wrong values
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right and Wrong values, plus formatted Max string
$ws.Range("B12").Value = 115
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "109.0/140"
